$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump version text in A5
$ws.Range("A5").Value = "(Version: 1.0.1)"

# Update "last tested with" text in A6 (keep trailing space as in source)
$ws.Range("A6").Value = "(Last tested with: ReportServer 4.0.0-6053) "

# Move the active selection from A4 to A5, matching the saved selection in the file
$ws.Range("A5").Select()
